# Denied file upload for objects with DELETE / DELETED status.
#
# The "Complaint - Participants cannot upload or replace files" /
# "Case File - Participants cannot version files" /
# "Task - Participants cannot version files" rules used to only fire when
# status == 'CLOSED' (and only stop firing when status != 'CLOSED'). Extend
# both guard expressions so the deny rule also applies when the object's
# status is DELETE or DELETED.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$notClosedOrDeleted = "status != 'CLOSED' && status != 'DELETE' && status != 'DELETED'"
$closedOrDeleted     = "status == 'CLOSED' || status == 'DELETE' || status == 'DELETED'"

# Complaint rule pair (rows 26-27)
$ws.Range("D26").Value = $notClosedOrDeleted
$ws.Range("D26").WrapText = $true

$ws.Range("D27").Value = $closedOrDeleted
$ws.Range("D27").WrapText = $true

# Case File rule pair (rows 41-42)
$ws.Range("D41").Value = $notClosedOrDeleted
$ws.Range("D41").WrapText = $true

$ws.Range("D42").Value = $closedOrDeleted
$ws.Range("D42").WrapText = $true

# Task rule pair (rows 50-51)
$ws.Range("D50").Value = $notClosedOrDeleted
$ws.Range("D50").WrapText = $true

$ws.Range("D51").Value = $closedOrDeleted
$ws.Range("D51").WrapText = $true

# Match the author's final cursor position/selection in the sheet.
$ws.Range("D51").Select()
